# Weekly update: insert two new rows of "Palta" price data at the top of the
# existing block (rows 323-324), pushing the previous rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 323; Excel copies formatting (incl. the
# date style on column D) from the row above, same as it would interactively.
$ws.Rows("323:324").Insert()

# New row 323: "Primera" quality, week of 2021-11-09 (serial 44509)
$ws.Cells.Item(323, 1).Value = 7
$ws.Cells.Item(323, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(323, 3).Value = "Ñuble"
$ws.Cells.Item(323, 4).Value = 44509
$ws.Cells.Item(323, 5).Value = 16
$ws.Cells.Item(323, 6).Value = "Fruta"
$ws.Cells.Item(323, 7).Value = 100106
$ws.Cells.Item(323, 8).Value = "Oleaginosos"
$ws.Cells.Item(323, 9).Value = 100106002
$ws.Cells.Item(323, 10).Value = "Palta"
$ws.Cells.Item(323, 11).Value = "Hass"
$ws.Cells.Item(323, 12).Value = "Primera"
$ws.Cells.Item(323, 13).Value = 120
$ws.Cells.Item(323, 14).Value = 2800
$ws.Cells.Item(323, 15).Value = 2900
$ws.Cells.Item(323, 16).Value = 2850
$ws.Cells.Item(323, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(323, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(323, 19).Value = 2850
$ws.Cells.Item(323, 20).Value = 1

# New row 324: "Segunda" quality, same week
$ws.Cells.Item(324, 1).Value = 7
$ws.Cells.Item(324, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(324, 3).Value = "Ñuble"
$ws.Cells.Item(324, 4).Value = 44509
$ws.Cells.Item(324, 5).Value = 16
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100106
$ws.Cells.Item(324, 8).Value = "Oleaginosos"
$ws.Cells.Item(324, 9).Value = 100106002
$ws.Cells.Item(324, 10).Value = "Palta"
$ws.Cells.Item(324, 11).Value = "Hass"
$ws.Cells.Item(324, 12).Value = "Segunda"
$ws.Cells.Item(324, 13).Value = 120
$ws.Cells.Item(324, 14).Value = 2500
$ws.Cells.Item(324, 15).Value = 2600
$ws.Cells.Item(324, 16).Value = 2550
$ws.Cells.Item(324, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(324, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(324, 19).Value = 2550
$ws.Cells.Item(324, 20).Value = 1
